# Correcoes Tela / Adicionamento de Dependencias
#
# - Remove the old "rede1/rede2/rede3" sample rows and replace the
#   network values with "pokerstars".
# - Row 2 and row 3 now hold the same tournament id (3540607900) and
#   the same network name ("pokerstars"), both written as plain,
#   unstyled cells.
# - The old trailing sample row (row 4) is removed entirely.
# - The stray formatted-but-empty cell that used to live at C3 is
#   removed, and a new empty, underlined placeholder cell is added at
#   E3 instead (this becomes the new active/selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last sample row ("3", "rede3") completely - shifts nothing
# else since it was the last row.
$ws.Rows.Item(4).Delete()

# Row 2: tournament id + network, no special formatting.
$ws.Range("A2").Value = 3540607900
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = "pokerstars"
$ws.Range("B2").ClearFormats()

# Row 3: same tournament id + network, no special formatting.
$ws.Range("A3").Value = 3540607900
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "pokerstars"
$ws.Range("B3").ClearFormats()

# Remove the old empty styled cell at C3.
$ws.Range("C3").Clear()

# Add the new empty, underlined cell at E3.
$ws.Range("E3").Font.Underline = $true

# Make E3 the active selection, matching the saved view state.
$ws.Range("E3").Select()
